$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 174
$ws.Range("F6").Value = 3194
$ws.Range("F7").Value = 2763
$ws.Range("F10").Value = 22
$ws.Range("F11").Value = 358
$ws.Range("F12").Value = 296
$ws.Range("F14").Value = 5774
$ws.Range("F16").Value = 1027
$ws.Range("F17").Value = 64
$ws.Range("F20").Value = 476
$ws.Range("F21").Value = 1255
$ws.Range("F23").Value = 7
$ws.Range("F24").Value = 1765
$ws.Range("F25").Value = 138
$ws.Range("F26").Value = 337

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 1171
$ws.Range("F9").Value = 66
$ws.Range("F13").Value = 636
$ws.Range("F24").Value = 293
$ws.Range("F25").Value = 4035
$ws.Range("F29").Value = 210
$ws.Range("F30").Value = 64

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 2569
$ws.Range("F6").Value = 1121
$ws.Range("F9").Value = 1450
$ws.Range("F11").Value = 113
$ws.Range("F13").Value = 555

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 2569
$ws.Range("F6").Value = 1121
$ws.Range("F7").Value = 1450
$ws.Range("F9").Value = 113
$ws.Range("F14").Value = 3194
$ws.Range("F15").Value = 2763
$ws.Range("F18").Value = 22
$ws.Range("F19").Value = 358
$ws.Range("F21").Value = 66
$ws.Range("F22").Value = 296
$ws.Range("F26").Value = 1027
$ws.Range("F27").Value = 636
$ws.Range("F28").Value = 64
$ws.Range("F31").Value = 476
$ws.Range("F39").Value = 1255
$ws.Range("F41").Value = 210
$ws.Range("F42").Value = 64
$ws.Range("F43").Value = 1768
$ws.Range("F46").Value = 138
$ws.Range("F47").Value = 337
